$d = $word.ActiveDocument

function Get-ParagraphByText([string]$searchText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Paragraph containing '$searchText' not found"
    }
    return $rng.Paragraphs(1)
}

function Set-ParagraphXml($paragraph, [string]$bodyXml) {
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $paragraph.Range.InsertXML($xml)
}

# --- 1) "See myModule.js file for example. Example for using module.exports." ---
# Split " Example for using module.exports." into three runs with proofErr markers
# around "module.exports".
$p1 = Get-ParagraphByText("Example for using module.exports.")
$body1 = '<w:p><w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>See myModule.js file for example.</w:t></w:r>' +
    '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> Example for using </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>module.exports</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>.</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $p1 $body1

# --- 2) "Custom events with the EventEmitter" heading, split + new paragraph ---
$p2 = Get-ParagraphByText("Custom events with the EventEmitter")
$body2 = '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Custom events with the </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>EventEmitter</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>' +
    '<w:p>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>EventEmitter</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> gives us the option to create and handle custom events.</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $p2 $body2

# --- 3) "Improve a module with EventEmiter" heading ---
$p3 = Get-ParagraphByText("Improve a module with EventEmiter")
$body3 = '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Improve a module with </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>EventEmiter</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'
Set-ParagraphXml $p3 $body3
